$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 156 - CESAMO San Manuel (0509)
$ws.Range("I156").Value = 9
$ws.Range("J156").Value = "0509"
$ws.Range("K156").Value = "San Manuel"
$ws.Range("M156").Value = 1
$ws.Range("N156").Value = "050901"
$ws.Range("O156").Value = "San Manuel"
$ws.Range("Q156").Value = "HND-0509"
$ws.Range("V156").Value = 15.329428999999999
$ws.Range("W156").Value = -87.921024000000003

# Row 157 - CESAMO Villanueva (0511)
$ws.Range("I157").Value = 11
$ws.Range("J157").Value = "0511"
$ws.Range("K157").Value = "Villanueva"
$ws.Range("M157").Value = 1
$ws.Range("N157").Value = "051101"
$ws.Range("O157").Value = "Villanueva"
$ws.Range("Q157").Value = "HND-0511"
$ws.Range("V157").Value = 15.312935
$ws.Range("W157").Value = -87.993703999999994

# Row 158 - CESAMO Potrerillos (0505)
$ws.Range("I158").Value = 5
$ws.Range("J158").Value = "0505"
$ws.Range("K158").Value = "Potrerillos"
$ws.Range("M158").Value = 1
$ws.Range("N158").Value = "050501"
$ws.Range("O158").Value = "Potrerillos"
$ws.Range("Q158").Value = "HND-0505"
$ws.Range("V158").Value = 15.228071
$ws.Range("W158").Value = -87.964016999999998

# Row 159 - CESAMO San Antonio de Cortés (0507)
$ws.Range("I159").Value = 7
$ws.Range("J159").Value = "0507"
$ws.Range("K159").Value = "San Antonio de Cortés"
$ws.Range("M159").Value = 1
$ws.Range("N159").Value = "050701"
$ws.Range("O159").Value = "San Antonio de Cortés"
$ws.Range("Q159").Value = "HND-0507"
$ws.Range("V159").Value = 15.114108
$ws.Range("W159").Value = -88.040538999999995

# Row 160 - CESAMO Santa Cruz de Yojoa (0510); also fill in department cells E/F/G
$ws.Range("E160").Value = 5
$ws.Range("F160").Value = "05"
$ws.Range("G160").Value = "Cortés"
$ws.Range("I160").Value = 10
$ws.Range("J160").Value = "0510"
$ws.Range("K160").Value = "Santa Cruz de Yojoa"
$ws.Range("M160").Value = 1
$ws.Range("N160").Value = "051001"
$ws.Range("O160").Value = "Santa Cruz de Yojoa"
$ws.Range("Q160").Value = "HND-0510"
$ws.Range("V160").Value = 14.97907
$ws.Range("W160").Value = -87.890918999999997

# Move selection to reflect the final edited cell (mirrors the saved view state)
$ws.Range("W162").Select()
